$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "username" and "email" values in row 7/8 from halosalsa1 -> halosalsa2
$ws.Range("B7").Value = "halosalsa2"
$ws.Range("B8").Value = "halosalsa2@gmail.com"
